$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 135.7046583333333
$ws.Range("H2").Value = 407.113975
$ws.Range("I2").Value = 0.2901853119378819
$ws.Range("J2").Value = 0.2901853119378819
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 98.946724
$ws.Range("N2").Value = 296.840172
$ws.Range("O2").Value = 0.2098009692989996
$ws.Range("P2").Value = 0.2098009692989996
$ws.Range("Q2").Value = 13427.53137362263
$ws.Range("R2").Value = 120847.7823626037
$ws.Range("S2").Value = 0.06088115972090018
$ws.Range("T2").Value = 0.06088115972090019
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 135.7046583333333
$ws.Range("H3").Value = 407.113975
$ws.Range("I3").Value = 0.2901853119378819
$ws.Range("J3").Value = 0.2901853119378819
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 163.0062356666667
$ws.Range("N3").Value = 489.018707
$ws.Range("O3").Value = 0.345629090707923
$ws.Range("P3").Value = 0.345629090707923
$ws.Range("Q3").Value = 22120.70551734781
$ws.Range("R3").Value = 199086.3496561303
$ws.Range("S3").Value = 0.1002964855018851
$ws.Range("T3").Value = 0.1002964855018851
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 135.7046583333333
$ws.Range("H4").Value = 407.113975
$ws.Range("I4").Value = 0.2901853119378819
$ws.Range("J4").Value = 0.2901853119378819
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 65.39610666666668
$ws.Range("N4").Value = 196.18832
$ws.Range("O4").Value = 0.1386621609326595
$ws.Range("P4").Value = 0.1386621609326595
$ws.Range("Q4").Value = 8874.556311530225
$ws.Range("R4").Value = 79871.006803772
$ws.Range("S4").Value = 0.04023772242422457
$ws.Range("T4").Value = 0.04023772242422458
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 135.7046583333333
$ws.Range("H5").Value = 407.113975
$ws.Range("I5").Value = 0.2901853119378819
$ws.Range("J5").Value = 0.2901853119378819
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 144.2727966666667
$ws.Range("N5").Value = 432.81839
$ws.Range("O5").Value = 0.3059077790604178
$ws.Range("P5").Value = 0.3059077790604179
$ws.Range("Q5").Value = 19578.49057844447
$ws.Range("R5").Value = 176206.4152060002
$ws.Range("S5").Value = 0.08876994429087201
$ws.Range("T5").Value = 0.08876994429087204
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 226.082006
$ws.Range("H6").Value = 678.246018
$ws.Range("I6").Value = 0.4834445496594812
$ws.Range("J6").Value = 0.4834445496594812
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 98.946724
$ws.Range("N6").Value = 296.840172
$ws.Range("O6").Value = 0.2098009692989996
$ws.Range("P6").Value = 0.2098009692989996
$ws.Range("Q6").Value = 22370.07384904835
$ws.Range("R6").Value = 201330.6646414351
$ws.Range("S6").Value = 0.1014271351208775
$ws.Range("T6").Value = 0.1014271351208775
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 226.082006
$ws.Range("H7").Value = 678.246018
$ws.Range("I7").Value = 0.4834445496594812
$ws.Range("J7").Value = 0.4834445496594812
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 163.0062356666667
$ws.Range("N7").Value = 489.018707
$ws.Range("O7").Value = 0.345629090707923
$ws.Range("P7").Value = 0.345629090707923
$ws.Range("Q7").Value = 36852.77675002875
$ws.Range("R7").Value = 331674.9907502587
$ws.Range("S7").Value = 0.1670925001065078
$ws.Range("T7").Value = 0.1670925001065079
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 226.082006
$ws.Range("H8").Value = 678.246018
$ws.Range("I8").Value = 0.4834445496594812
$ws.Range("J8").Value = 0.4834445496594812
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 65.39610666666668
$ws.Range("N8").Value = 196.18832
$ws.Range("O8").Value = 0.1386621609326595
$ws.Range("P8").Value = 0.1386621609326595
$ws.Range("Q8").Value = 14784.88297978998
$ws.Range("R8").Value = 133063.9468181098
$ws.Range("S8").Value = 0.06703546594690006
$ws.Range("T8").Value = 0.06703546594690007
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 226.082006
$ws.Range("H9").Value = 678.246018
$ws.Range("I9").Value = 0.4834445496594812
$ws.Range("J9").Value = 0.4834445496594812
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 144.2727966666667
$ws.Range("N9").Value = 432.81839
$ws.Range("O9").Value = 0.3059077790604178
$ws.Range("P9").Value = 0.3059077790604179
$ws.Range("Q9").Value = 32617.48328163011
$ws.Range("R9").Value = 293557.3495346711
$ws.Range("S9").Value = 0.1478894484851958
$ws.Range("T9").Value = 0.1478894484851958
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.1600446666666666
$ws.Range("H10").Value = 0.4801339999999999
$ws.Range("I10").Value = 0.0003422329939962955
$ws.Range("J10").Value = 0.0003422329939962955
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 98.946724
$ws.Range("N10").Value = 296.840172
$ws.Range("O10").Value = 0.2098009692989996
$ws.Range("P10").Value = 0.2098009692989996
$ws.Range("Q10").Value = 15.83589546033866
$ws.Range("R10").Value = 142.523059143048
$ws.Range("S10").Value = 0.0000718008138665215
$ws.Range("T10").Value = 0.0000718008138665215
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.1600446666666666
$ws.Range("H11").Value = 0.4801339999999999
$ws.Range("I11").Value = 0.0003422329939962955
$ws.Range("J11").Value = 0.0003422329939962955
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 163.0062356666667
$ws.Range("N11").Value = 489.018707
$ws.Range("O11").Value = 0.345629090707923
$ws.Range("P11").Value = 0.345629090707923
$ws.Range("Q11").Value = 26.08827865185977
$ws.Range("R11").Value = 234.794507866738
$ws.Range("S11").Value = 0.0001182856785251897
$ws.Range("T11").Value = 0.0001182856785251897
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.1600446666666666
$ws.Range("H12").Value = 0.4801339999999999
$ws.Range("I12").Value = 0.0003422329939962955
$ws.Range("J12").Value = 0.0003422329939962955
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 65.39610666666668
$ws.Range("N12").Value = 196.18832
$ws.Range("O12").Value = 0.1386621609326595
$ws.Range("P12").Value = 0.1386621609326595
$ws.Range("Q12").Value = 10.46629809276445
$ws.Range("R12").Value = 94.19668283488001
$ws.Range("S12").Value = 0.0000474547664899802
$ws.Range("T12").Value = 0.00004745476648998021
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.1600446666666666
$ws.Range("H13").Value = 0.4801339999999999
$ws.Range("I13").Value = 0.0003422329939962955
$ws.Range("J13").Value = 0.0003422329939962955
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 144.2727966666667
$ws.Range("N13").Value = 432.81839
$ws.Range("O13").Value = 0.3059077790604178
$ws.Range("P13").Value = 0.3059077790604179
$ws.Range("Q13").Value = 23.09009165158444
$ws.Range("R13").Value = 207.81082486426
$ws.Range("S13").Value = 0.0001046917351146041
$ws.Range("T13").Value = 0.0001046917351146041
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 105.7015583333333
$ws.Range("H14").Value = 317.104675
$ws.Range("I14").Value = 0.2260279054086406
$ws.Range("J14").Value = 0.2260279054086406
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 98.946724
$ws.Range("N14").Value = 296.840172
$ws.Range("O14").Value = 0.2098009692989996
$ws.Range("P14").Value = 0.2098009692989996
$ws.Range("Q14").Value = 10458.82291877823
$ws.Range("R14").Value = 94129.40626900409
$ws.Range("S14").Value = 0.04742087364335539
$ws.Range("T14").Value = 0.04742087364335539
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 105.7015583333333
$ws.Range("H15").Value = 317.104675
$ws.Range("I15").Value = 0.2260279054086406
$ws.Range("J15").Value = 0.2260279054086406
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 163.0062356666667
$ws.Range("N15").Value = 489.018707
$ws.Range("O15").Value = 0.345629090707923
$ws.Range("P15").Value = 0.3456290907079231
$ws.Range("Q15").Value = 17230.01312801725
$ws.Range("R15").Value = 155070.1181521552
$ws.Range("S15").Value = 0.07812181942100488
$ws.Range("T15").Value = 0.07812181942100489
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 105.7015583333333
$ws.Range("H16").Value = 317.104675
$ws.Range("I16").Value = 0.2260279054086406
$ws.Range("J16").Value = 0.2260279054086406
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 65.39610666666668
$ws.Range("N16").Value = 196.18832
$ws.Range("O16").Value = 0.1386621609326595
$ws.Range("P16").Value = 0.1386621609326595
$ws.Range("Q16").Value = 6912.470383599557
$ws.Range("R16").Value = 62212.23345239601
$ws.Range("S16").Value = 0.03134151779504485
$ws.Range("T16").Value = 0.03134151779504486
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 105.7015583333333
$ws.Range("H17").Value = 317.104675
$ws.Range("I17").Value = 0.2260279054086406
$ws.Range("J17").Value = 0.2260279054086406
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 144.2727966666667
$ws.Range("N17").Value = 432.81839
$ws.Range("O17").Value = 0.3059077790604178
$ws.Range("P17").Value = 0.3059077790604179
$ws.Range("Q17").Value = 15249.8594327748
$ws.Range("R17").Value = 137248.7348949733
$ws.Range("S17").Value = 0.06914369454923544
$ws.Range("T17").Value = 0.06914369454923545
